# Applies the cryptos-list price/volume refresh described in the commit.
# Column D ("Price") cells that are plain decimal numbers (e.g. "227.64") are
# forced to Text format first so Excel keeps the exact display string instead
# of silently re-interpreting them as floating-point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.461.35"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "1.803.42"
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.64"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.575"
$ws.Range("E6").Value = "  +3.24%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "36.46"
$ws.Range("E8").Value = "  +5.31%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.298"
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0691"
$ws.Range("E10").Value = "  -0.41%  "
$ws.Range("E11").Value = "  +1.40%  "
$ws.Range("D12").Value = "2.062.99"
$ws.Range("E12").Value = "  -0.99%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.60"
$ws.Range("E13").Value = "  +2.16%  "
$ws.Range("D14").Value = "1.815.46"
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.646"
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("E16").Value = "  +3.91%  "
$ws.Range("D17").Value = "34.418.88"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.16"
$ws.Range("E18").Value = "  +1.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.76"
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("D20").Value = "0.0₃0791"
$ws.Range("E20").Value = "  -1.41%  "
$ws.Range("E21").Value = "  +0.45%  "
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.18"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("E24").Value = "  +4.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "172.32"
$ws.Range("E25").Value = "  +2.33%  "
$ws.Range("E26").Value = "  +8.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.98"
$ws.Range("E27").Value = "  +1.30%  "
$ws.Range("E28").Value = "  +1.62%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.02"
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.85"
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0529"
$ws.Range("E32").Value = "  -0.36%  "
$ws.Range("E33").Value = "  -0.50%  "
$ws.Range("E34").Value = "  -2.22%  "
$ws.Range("D35").Value = "1.396.51"
$ws.Range("E35").Value = "  -1.33%  "
$ws.Range("E36").Value = "  -0.80%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.44"
$ws.Range("E37").Value = "  -6.63%  "
$ws.Range("E38").Value = "  -0.50%  "
$ws.Range("E39").Value = "  -0.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "82.62"
$ws.Range("E40").Value = "  -4.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.960"
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("E42").Value = "  -1.16%  "
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("E44").Value = "  +7.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.36"
$ws.Range("E45").Value = "  -4.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.03"
$ws.Range("E46").Value = "  -0.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0503"
$ws.Range("E47").Value = "  -4.14%  "
$ws.Range("D48").Value = "1.964.02"
$ws.Range("E48").Value = "  -0.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "104.43"
$ws.Range("E49").Value = "  -1.45%  "
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("D51").Value = "0.0₆0128"
$ws.Range("E51").Value = "  -0.59%  "
